$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    # Force the cell to remain a text value (rather than letting Excel's
    # General-format auto-detection turn a numeric-looking string into a
    # real number), then restore the default "Normal" style so no
    # formatting side effects are introduced.
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

# Column D ("Price") updates - values look numeric, so they must be
# written through Set-TextValue to keep them stored as text.
Set-TextValue "D3"  "23.94"
Set-TextValue "D4"  "5.204"
Set-TextValue "D5"  "0.05732"
Set-TextValue "D6"  "6.488"
Set-TextValue "D7"  "3.165"
Set-TextValue "D8"  "0.8146"
Set-TextValue "D9"  "0.8699"
Set-TextValue "D10" "0.1368"
Set-TextValue "D11" "0.06937"
Set-TextValue "D12" "0.03171"
Set-TextValue "D14" "0.09322"
Set-TextValue "D15" "3.862"
Set-TextValue "D16" "0.001524"
Set-TextValue "D18" "0.0005973"
Set-TextValue "D19" "0.006163"
Set-TextValue "D20" "0.001239"
Set-TextValue "D22" "0.00008504"
Set-TextValue "D24" "2.163"
Set-TextValue "D27" "0.0002332"
Set-TextValue "D41" "0.006373"
Set-TextValue "D44" "0.008104"
Set-TextValue "D45" "0.00005474"
Set-TextValue "D47" "0.4541"
Set-TextValue "D48" "0.002522"

# Column E ("Volume(1h)") updates - plain alphanumeric text, no special
# handling required since Excel won't mistake these for numbers.
$ws.Range("E18").Value2 = "17OneONE"
$ws.Range("E41").Value2 = "40KickTokenKICKBestin24h"
$ws.Range("E47").Value2 = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("E48").Value2 = "47BOLOBOLO"
